$wb = $excel.ActiveWorkbook

# Both the "展览" and "全部类型" sheets carry the same table of events;
# update the "想去人数" (column F) counts on each.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 125
    $ws.Range("F8").Value = 11586
    $ws.Range("F10").Value = 93
    $ws.Range("F11").Value = 458
    $ws.Range("F12").Value = 374
    $ws.Range("F16").Value = 13158
}
